$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H17").Value = 19532.053
$ws.Range("J17").Value = 20517.166
$ws.Range("L17").Value = 61551.49800000001
$ws.Range("N17").Value = -61887.49800000001

$ws.Range("H70").Value = 4174.25
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 4174.25
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 12522.75
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -13062.75

$ws.Range("H73").Value = 4174.25
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 4174.25
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 12522.75
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -14394.75

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H125").Value = 2491.75
$ws.Range("J125").Value = 967
$ws.Range("L125").Value = 8703
$ws.Range("N125").Value = -13623

$ws.Range("H135").Value = 1362.65
$ws.Range("I135").Value = 791.8889
$ws.Range("K135").Value = 7127.0001
$ws.Range("M135").Value = -4592.0001

$ws.Range("H137").Value = 45580.05
$ws.Range("I137").Value = 65610.16
$ws.Range("J137").Value = 2181.5
$ws.Range("K137").Value = 196830.48
$ws.Range("L137").Value = 6544.5
$ws.Range("M137").Value = -194280.48
$ws.Range("N137").Value = -11644.5

$ws.Range("H140").Value = 171915.67
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 171915.67
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 171915.67
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -182275.67

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H21").Value = 12390
$ws.Range("I21").Value = 3983.3333
$ws.Range("K21").Value = 3983.3333
$ws.Range("M21").Value = -3609.3333

$ws.Range("H74").Value = 2657.3333
$ws.Range("I74").Value = 2665.077
$ws.Range("K74").Value = 2665.077
$ws.Range("M74").Value = -1791.077

$ws.Range("H77").Value = 2657.3333
$ws.Range("I77").Value = 2665.077
$ws.Range("K77").Value = 13325.385
$ws.Range("M77").Value = -8957.385000000002

$ws.Range("H97").Value = 977.7727
$ws.Range("I97").Value = 605.8421
$ws.Range("J97").Value = 3333.3333
$ws.Range("K97").Value = 605.8421
$ws.Range("L97").Value = 3333.3333
$ws.Range("M97").Value = -109.8421
$ws.Range("N97").Value = -4325.3333

$ws.Range("H107").Value = 38999.5
$ws.Range("J107").Value = 38999.5
$ws.Range("L107").Value = 38999.5
$ws.Range("N107").Value = -46679.5

$ws.Range("H110").Value = 1971.7
$ws.Range("I110").Value = 1131
$ws.Range("J110").Value = 3933.3333
$ws.Range("K110").Value = 1131
$ws.Range("L110").Value = 3933.3333
$ws.Range("M110").Value = 914
$ws.Range("N110").Value = -8023.3333

$ws.Range("H128").Value = 99889.5
$ws.Range("J128").Value = 99889.5
$ws.Range("L128").Value = 99889.5
$ws.Range("N128").Value = -109849.5

$ws.Range("H131").Value = 94993.8
$ws.Range("I131").Value = 94993
$ws.Range("K131").Value = 94993
$ws.Range("M131").Value = -89953

$ws.Range("H133").Value = 117951.375
$ws.Range("J133").Value = 118928.57
$ws.Range("L133").Value = 118928.57
$ws.Range("N133").Value = -123988.57

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H107").Value = 1201.1818
$ws.Range("I107").Value = 596.375
$ws.Range("J107").Value = 2814
$ws.Range("K107").Value = 596.375
$ws.Range("L107").Value = 2814
$ws.Range("M107").Value = 1323.625
$ws.Range("N107").Value = -6654

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H75").Value = 90424.71000000001
$ws.Range("J75").Value = 102795
$ws.Range("L75").Value = 102795
$ws.Range("N75").Value = -104791

$ws.Range("H78").Value = 90424.71000000001
$ws.Range("J78").Value = 102795
$ws.Range("L78").Value = 308385
$ws.Range("N78").Value = -318369

$ws.Range("H105").Value = 4000
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()

$ws.Range("H134").Value = 2155.3513
$ws.Range("I134").Value = 1621.0646
$ws.Range("K134").Value = 4863.1938
$ws.Range("M134").Value = -2328.1938

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H113").Value = 1496.5769
$ws.Range("J113").Value = 1444.5264
$ws.Range("L113").Value = 4333.5792
$ws.Range("N113").Value = -8673.5792

$ws.Range("H131").Value = 3204.5386
$ws.Range("I131").Value = 1874.75
$ws.Range("K131").Value = 5624.25
$ws.Range("M131").Value = -584.25

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H75").Value = 111969
$ws.Range("J75").Value = 111969
$ws.Range("L75").Value = 111969
$ws.Range("N75").Value = -113717

$ws.Range("H78").Value = 111969
$ws.Range("J78").Value = 111969
$ws.Range("L78").Value = 335907
$ws.Range("N78").Value = -344643

$ws.Range("H128").Value = 152700
$ws.Range("J128").Value = 152700
$ws.Range("L128").Value = 152700
$ws.Range("N128").Value = -162660

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H46").Value = 3529.3572
$ws.Range("I46").Value = 1337.4166
$ws.Range("K46").Value = 1337.4166
$ws.Range("M46").Value = -1149.4166

$ws.Range("H61").Value = 10833.333
$ws.Range("I61").Value = 10833.333
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 10833.333
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -10631.333
$ws.Range("N61").ClearContents()

$ws.Range("H113").Value = 10833.333
$ws.Range("I113").Value = 10833.333
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 10833.333
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -8663.333000000001
$ws.Range("N113").ClearContents()

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws.Range("H139").Value = 95233.336
$ws.Range("J139").Value = 95233.336
$ws.Range("L139").Value = 95233.336
$ws.Range("N139").Value = -105513.336

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H31").Value = 40876.25
$ws.Range("J31").Value = 59019
$ws.Range("L31").Value = 59019
$ws.Range("N31").Value = -59715

$ws.Range("H81").Value = 3850.5
$ws.Range("J81").Value = 4346.6924
$ws.Range("L81").Value = 8693.3848
$ws.Range("N81").Value = -10815.3848

$ws.Range("H84").Value = 3850.5
$ws.Range("J84").Value = 4346.6924
$ws.Range("L84").Value = 43466.924
$ws.Range("N84").Value = -54074.924

$ws.Range("H113").Value = 533.5
$ws.Range("I113").Value = 469.6
$ws.Range("K113").Value = 1408.8
$ws.Range("M113").Value = 761.1999999999998

$ws.Range("H126").Value = 9399
$ws.Range("J126").Value = 9498.75
$ws.Range("L126").Value = 28496.25
$ws.Range("N126").Value = -33436.25

$ws.Range("H128").Value = 165854.75
$ws.Range("J128").Value = 165854.75
$ws.Range("L128").Value = 165854.75
$ws.Range("N128").Value = -175814.75

$ws.Range("H130").Value = 122981
$ws.Range("J130").Value = 122981
$ws.Range("L130").Value = 122981
$ws.Range("N130").Value = -133021
